# Rename the "localtestdir" argument-name entry to "locdir"
# (commit message: "locdir used in website").
#
# The text lives in a run inside one paragraph of the
# "Rectangle: Rounded Corners 98" auto-shape, which stacks several
# function-argument names as separate paragraphs. We search every shape on
# slide 1 (recursing into groups, defensively) for a run whose text is
# exactly "localtestdir" and rename it, leaving every other paragraph/run
# untouched.
#
# Note: Paragraph.Text includes a trailing "`r" paragraph-mark character
# (PowerPoint COM convention), so comparisons against paragraph text need
# to strip that before comparing; Run.Text does not include it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$script:found = $false

function Update-Shape($shp) {
    if ($shp.Type -eq 6) {
        # msoGroup -> recurse into its members
        $giCount = $shp.GroupItems.Count
        for ($gi = 1; $gi -le $giCount; $gi++) {
            Update-Shape($shp.GroupItems.Item($gi))
        }
        return
    }

    if ($shp.HasTextFrame -eq 0) { return }
    $tf = $shp.TextFrame
    if ($tf.HasText -eq 0) { return }

    $tr = $tf.TextRange
    $parasCount = $tr.Paragraphs().Count
    for ($pi = 1; $pi -le $parasCount; $pi++) {
        $para = $tr.Paragraphs($pi)
        $paraText = $para.Text.TrimEnd("`r")
        if ($paraText -ne "localtestdir") { continue }

        $runsCount = $para.Runs().Count
        for ($ri = 1; $ri -le $runsCount; $ri++) {
            $run = $para.Runs($ri)
            if ($run.Text -eq "localtestdir") {
                $run.Text = "locdir"
                $script:found = $true
            }
        }
    }
}

for ($si = 1; $si -le $s.Shapes.Count; $si++) {
    Update-Shape($s.Shapes.Item($si))
}

if (-not $script:found) {
    throw "Could not find run with text 'localtestdir' to update"
}
